$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Planilha1 -> FP)
$ws.Name = "FP"

# New "Esforço / Prazo / Custo" summary block (rows 30-31) and
# "Salario médio analista" row (row 33).
# Values are written in the same order the shared-string table expects
# them to be interned (Esforço, 5 FPs/mês, Prazo, 3,4 meses, Custo, ...).
$ws.Range("A30").Value = "Esforço"
$ws.Range("A31").Value = "5 FPs/mês"

$ws.Range("B30").Value = "Prazo "
$ws.Range("B31").Value = "3,4 meses"

$ws.Range("C30").Value = "Custo"
$ws.Range("C31").Value = 15688.08
$ws.Range("C31").NumberFormat = """R$""\ #,##0.00;[Red]\-""R$""\ #,##0.00"

$ws.Range("A33").Value = "Salario médio analista "
$ws.Range("B33").Value = 3922.02
$ws.Range("B33").NumberFormat = """R$""\ #,##0.00;[Red]\-""R$""\ #,##0.00"

# Move the active selection the way the authored workbook left it.
$ws.Range("C34").Select() | Out-Null
